$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: simple text replacements (row indices below are the
# original 1-based row numbers; none of these operations change the
# row count, so the indices stay valid throughout this step) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text  = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text  = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text  = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text  = "1788"
$t.Rows.Item(6).Cells.Item(1).Range.Text  = "0.00019"
$t.Rows.Item(7).Cells.Item(1).Range.Text  = "0.00020"
$t.Rows.Item(8).Cells.Item(1).Range.Text  = "0.00022"
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.22655"

# Last three rows of the table used to hold one run per value separated
# by tabs; they collapse down to a single plain value each.
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.94"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.23"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "378"

# --- Step 2: drop the three now-superfluous rows that used to sit
# right after the (old) "0.00002" row (rows 9, 10, 11 of the original
# table: "0.00011", "0.00011", "0.00012"). Deleting index 9 three times
# removes all three because each delete shifts the following rows up
# by one. ---
$t.Rows.Item(9).Delete()
$t.Rows.Item(9).Delete()
$t.Rows.Item(9).Delete()

# --- Step 3: insert three brand new rows right before the row that
# still holds "0.00003" (originally row 5, untouched and therefore
# still row 5 after the deletes above, since the deletes happened
# after it). Rows.Add(ref) always inserts directly above "ref", so
# add them in reverse order to land in the desired final order:
# 0.00001, 0.00256, 0.00011. ---
$refRow = $t.Rows.Item(5)

$r3 = $t.Rows.Add($refRow)
$r3.Cells.Item(1).Range.Text = "0.00011"

$r2 = $t.Rows.Add($refRow)
$r2.Cells.Item(1).Range.Text = "0.00256"

$r1 = $t.Rows.Add($refRow)
$r1.Cells.Item(1).Range.Text = "0.00001"

Write-Output ("Final row count: " + $t.Rows.Count)
